$d = $word.ActiveDocument

function Split-PinDeclaration {
    param(
        [string]$OldText,
        [string[]]$NewParts
    )

    # Locate the existing run's text in the document.
    $find = $d.Content
    $found = $find.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $OldText"
    }

    $ranges = New-Object System.Collections.ArrayList

    # Replace the original run's text with the first new chunk.
    $r = $d.Range($find.Start, $find.End)
    $r.Text = $NewParts[0]
    [void]$ranges.Add($r)

    # Insert the remaining chunks as their own ranges, right after the previous one.
    $pos = $r.End
    for ($i = 1; $i -lt $NewParts.Length; $i++) {
        $ins = $d.Range($pos, $pos)
        $ins.InsertAfter($NewParts[$i])
        [void]$ranges.Add($ins)
        $pos = $ins.End
    }

    # Toggling a character attribute on/off after all the text is in place forces
    # each chunk to materialize as its own run (with identical formatting to the
    # original), instead of being re-merged into one run on save.
    foreach ($rng in $ranges) {
        $rng.Bold = 1
        $rng.Bold = 0
    }
}

Split-PinDeclaration " trigger1 = 2; //Trigger pin of 1st " @(" trigger1 = 11", "; //Trigger pin ", "of 1st ")
Split-PinDeclaration " echo1 = 3; //Echo pin of 1st " @(" echo1 = 10", "; //Echo pin of ", "1st ")
Split-PinDeclaration " trigger2 = 4; //Trigger pin of 2nd " @(" trigger2 = 6", "; //Trigger pin of 2nd ")

Write-Output "done"
